# Actualización de datos obtenidos el 6 de abril de 2016
#
# The "numero-habitaciones" mapping sheet pairs a numeric code (column A)
# with its canonical IAEST URI (column B, e.g. ".../numero-habitaciones/8093").
# A handful of rows had the wrong code/URI attached to them; this swaps/
# rotates the (A, B) pair of each affected row so every row is internally
# consistent again (A's number matches the trailing segment of B's URI).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseUri = "http://opendata.aragon.es/kos/iaest/numero-habitaciones/"

# row -> corrected numero-habitaciones code
$fixes = @{
    1   = 8093
    2   = 16190
    14  = 5
    15  = 1586
    16  = 10308
    40  = 1488
    41  = 1125
    47  = 827
    48  = 948
    49  = 134391
    100 = 754
    101 = 633
    102 = 2025
    103 = 1177
    104 = 3113
    115 = 887
    116 = 3487
    285 = 363
    286 = 242
    287 = 121
    288 = 484
    321 = 131
    322 = 252
    323 = 373
    349 = 39
    350 = 21099
    409 = 58
    410 = 5198
    514 = 1528
    515 = 17557
    524 = 2766
    525 = 5917
}

foreach ($row in $fixes.Keys) {
    $code = $fixes[$row]
    $ws.Cells.Item($row, 1).Value = $code
    $ws.Cells.Item($row, 2).Value = "$baseUri$code"
}
